$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Fix trailing space on paragraph 25 ("15. ... SeedResidentDemo.bat as administrator ")
# ------------------------------------------------------------------
$p25 = $d.Paragraphs.Item(25)
$r = $p25.Range.Duplicate
$r.Find.Execute(" as administrator ", $false, $false, $false, $false, $false, $true, 1, $false, " as administrator", 2) | Out-Null

# ------------------------------------------------------------------
# Step 2: Renumber paragraph 25: "15." -> "16."
# ------------------------------------------------------------------
$p25 = $d.Paragraphs.Item(25)
$r = $p25.Range.Duplicate
$r.Find.Execute("15. Run the script ", $false, $false, $false, $false, $false, $true, 1, $false, "16. Run the script ", 2) | Out-Null

# ------------------------------------------------------------------
# Step 3: Renumber paragraph 24: "14." -> "15."
# ------------------------------------------------------------------
$p24 = $d.Paragraphs.Item(24)
$r = $p24.Range.Duplicate
$r.Find.Execute("14. Run the script ", $false, $false, $false, $false, $false, $true, 1, $false, "15. Run the script ", 2) | Out-Null

# ------------------------------------------------------------------
# Step 4: Duplicate paragraph 23 ("13. Run the script C:\Deployments\Install\7_CreateTables.bat as administrator")
#         to a new paragraph right after it, which will become the new "14." step.
# ------------------------------------------------------------------
$p23 = $d.Paragraphs.Item(23)
$src = $p23.Range.Duplicate
$src.MoveEnd(1, -1) | Out-Null
$src.Copy()

$dst = $p23.Range.Duplicate
$dst.Collapse(0)
$dst.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(24)
$dstRange = $newPara.Range.Duplicate
$dstRange.Collapse(1)
$dstRange.Paste()

# Renumber the newly duplicated paragraph (still "13.") to "14."
$newPara = $d.Paragraphs.Item(24)
$r = $newPara.Range.Duplicate
$r.Find.Execute("13. Run the script ", $false, $false, $false, $false, $false, $true, 1, $false, "14. Run the script ", 2) | Out-Null

# ------------------------------------------------------------------
# Step 5: Rewrite paragraph 23 content to the new "uninstall" step text.
#   Before: "13. Run the script " + BOLD("C:\Deployments\Install\7_CreateTables.bat ") + "as administrator"
#   After:  "13. Run the script " + BOLD("C:\Deployments\Install\Utility\UninstallServices.bat") + " as administrator"
# ------------------------------------------------------------------
$p23 = $d.Paragraphs.Item(23)
$r = $p23.Range.Duplicate
$r.Find.Execute("C:\Deployments\Install\7_CreateTables.bat ", $false, $false, $false, $false, $false, $true, 1, $false, "C:\Deployments\Install\Utility\UninstallServices.bat", 2) | Out-Null

$p23 = $d.Paragraphs.Item(23)
$r = $p23.Range.Duplicate
$r.Find.Execute("as administrator", $false, $false, $false, $false, $false, $true, 1, $false, " as administrator", 2) | Out-Null

# ------------------------------------------------------------------
# Step 6: Append a brand new paragraph 17 at the end of the document:
#   "17. Run the script " + BOLD("C:\Deployments\Install\9_InstallServices.bat") + " as administrator"
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIndex)
$src2 = $pLast.Range.Duplicate
$src2.MoveEnd(1, -1) | Out-Null
$src2.Copy()

$dst2 = $pLast.Range.Duplicate
$dst2.Collapse(0)
$dst2.InsertParagraphAfter()

$newLastIndex = $d.Paragraphs.Count
$newPara2 = $d.Paragraphs.Item($newLastIndex)
$dstRange2 = $newPara2.Range.Duplicate
$dstRange2.Collapse(1)
$dstRange2.Paste()

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $newPara2.Range.Duplicate
$r.Find.Execute("16. Run the script ", $false, $false, $false, $false, $false, $true, 1, $false, "17. Run the script ", 2) | Out-Null

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $newPara2.Range.Duplicate
$r.Find.Execute("C:\Deployments\Install\Utility\SeedResidentDemo.bat", $false, $false, $false, $false, $false, $true, 1, $false, "C:\Deployments\Install\9_InstallServices.bat", 2) | Out-Null

Write-Host "Done. Total paragraphs:" $d.Paragraphs.Count
